$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11: Colaborador_id, Colaborador_nome, Departamento, Motivo_da_ausência, Horas_de_ausência, Data_da_ausência, Salário
$data = @(
    @(38823, "Pedro Miguel da Cruz", "Engenharia", "Problemas pessoais", 3, 45095, 6267.81),
    @(21520, "Gabrielly Moraes", "Atendimento ao Cliente", "Doença", 1, 45086, 3418),
    @(1713, "Lorenzo Silva", "Atendimento ao Cliente", "Consulta médica", 7, 45087, 5562.29),
    @(47308, "Dr. Vinicius Rodrigues", "Financeiro", "Consulta médica", 6, 45093, 6332.18),
    @(3087, "Murilo da Cunha", "TI", "Outros", 1, 45091, 5846.04),
    @(75021, "Srta. Sophia Ribeiro", "Marketing", "Problemas pessoais", 6, 45091, 11324.37),
    @(32589, "Dr. João Gabriel Teixeira", "P&D", "Doença", 8, 45086, 4079.4),
    @(86687, "Isabella das Neves", "TI", "Outros", 3, 45105, 8606.48),
    @(80543, "Dr. Bryan Cunha", "Engenharia", "Outros", 6, 45082, 11036.85),
    @(70264, "Sr. Levi Novaes", "Marketing", "Consulta médica", 8, 45078, 3543.44)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}

$wb.Save()
